# Refresh the cryptocurrency price / 1h-volume table with the latest
# values pulled by the scraper job. Row 23/24 (RenderToken <-> PancakeSwap)
# and row 31/32 (Bittensor <-> EthereumClassic) also swapped rank order,
# so their whole rows (coin name, link, price, volume) are rewritten.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.670.73"
$ws.Range("E2").Value = "  -2.64%  "
$ws.Range("D3").Value = "'3.185.04"
$ws.Range("E3").Value = "  -4.74%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'527.36"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").Value = "'172.67"
$ws.Range("E6").Value = "  -6.77%  "
$ws.Range("D7").Value = "'0.594"
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'3.181.22"
$ws.Range("E9").Value = "  -4.62%  "
$ws.Range("D10").Value = "'0.605"
$ws.Range("E10").Value = "  -2.22%  "
$ws.Range("D11").Value = "'53.21"
$ws.Range("E11").Value = "  -7.36%  "
$ws.Range("D12").Value = "'0.132"
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").Value = "'0.0000251"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").Value = "'9.06"
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").Value = "'3.687.10"
$ws.Range("E15").Value = "  -4.92%  "
$ws.Range("E16").Value = "  -3.22%  "
$ws.Range("D17").Value = "'3.169.59"
$ws.Range("E17").Value = "  -5.20%  "
$ws.Range("D18").Value = "'17.21"
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").Value = "'62.491.00"
$ws.Range("E19").Value = "  -2.50%  "
$ws.Range("D20").Value = "'11.00"
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("E21").Value = "  +1.38%  "
$ws.Range("D22").Value = "'365.42"
$ws.Range("E22").Value = "  -1.67%  "
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").Value = "'3.76"
$ws.Range("E23").Value = "  +2.10%  "
$ws.Range("B24").Value = "RenderToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D24").Value = "'11.22"
$ws.Range("E24").Value = "  +5.16%  "
$ws.Range("D25").Value = "'81.04"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("D26").Value = "'3.89"
$ws.Range("E26").Value = "  +4.29%  "
$ws.Range("D27").Value = "'6.11"
$ws.Range("E27").Value = "  +2.74%  "
$ws.Range("D28").Value = "'2.64"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("D29").Value = "'11.33"
$ws.Range("E29").Value = "  +1.60%  "
$ws.Range("D30").Value = "'8.18"
$ws.Range("E30").Value = "  -1.16%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'28.31"
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'639.41"
$ws.Range("E32").Value = "  -1.86%  "
$ws.Range("D33").Value = "'6.46"
$ws.Range("E33").Value = "  -2.82%  "
$ws.Range("D34").Value = "'11.31"
$ws.Range("E34").Value = "  +2.85%  "
$ws.Range("E35").Value = "  +2.27%  "
$ws.Range("D36").Value = "'56.35"
$ws.Range("E36").Value = "  -5.34%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").Value = "'37.00"
$ws.Range("E38").Value = "  +3.39%  "
$ws.Range("D39").Value = "'0.374"
$ws.Range("E39").Value = "  +0.78%  "
$ws.Range("D40").Value = "'0.996"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("D41").Value = "'0.0₃0706"
$ws.Range("E41").Value = "  +16.07%  "
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("D43").Value = "'2.876.76"
$ws.Range("E43").Value = "  +4.26%  "
$ws.Range("D44").Value = "'2.52"
$ws.Range("E44").Value = "  +10.47%  "
$ws.Range("D45").Value = "'2.93"
$ws.Range("E45").Value = "  +13.10%  "
$ws.Range("D46").Value = "'2.64"
$ws.Range("E46").Value = "  +2.64%  "
$ws.Range("D47").Value = "'0.0392"
$ws.Range("E47").Value = "  +2.57%  "
$ws.Range("D48").Value = "'2.58"
$ws.Range("E48").Value = "  -4.97%  "
$ws.Range("D49").Value = "'2.99"
$ws.Range("E49").Value = "  +8.01%  "
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("D51").Value = "'134.16"
$ws.Range("E51").Value = "  -0.51%  "

Write-Host "cryptos.xlsx refreshed: 102 cells updated"
